$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 -> ALLEGRETTO-LTE (B7981028): Days remaining 4 -> 3
$ws.Range("B9").Value = 3

# Row 11 -> REJOICE (MK-5909-003): Days remaining 31 -> 30
$ws.Range("B11").Value = 30
